$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: aggregate the DK1 entry into the combined DK entry
$ws.Range("A2").Value = "c_DK"
$ws.Range("B2").Value = "DK"
$ws.Range("C2").Value = "c_DK_Central"
$ws.Range("D2").Value = "DK"

# Delete the old row 3 (c_DK2 / DK2 / c_DK1_LargeDecentral / DK1), shifting
# the rows below up. This brings the already-blank A/B cells from the old
# row 4 into row 3 along with its C/D pair.
$ws.Range("A3:D3").Delete(-4162)  # xlShiftUp

# Row 3's C/D now become the single aggregated "Decentral" mapping entry
$ws.Range("C3").Value = "c_DK_Decentral"
$ws.Range("D3").Value = "DK"

# The remaining old rows (now rows 4-6, formerly the DK2 breakdown rows)
# are no longer needed now that everything rolls up into DK
$ws.Range("A4:D6").Delete(-4162)  # xlShiftUp
